# Auto-generated edit script: updates cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.221.62'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '3.269.13'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '531.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.73'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.595'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.04%  '
$ws.Range('D8').Value = '3.266.89'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.608'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.19'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -7.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.134'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000256'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.24'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.65%  '
$ws.Range('D15').Value = '3.805.11'
$ws.Range('E15').Value = '  +1.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.117'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.12%  '
$ws.Range('D17').Value = '3.282.19'
$ws.Range('E17').Value = '  +0.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.34'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('D19').Value = '63.218.91'
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.16'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.963'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '368.66'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.81%  '
$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.13'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +9.67%  '
$ws.Range('B24').Value = 'RenderToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.32'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.76'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '81.12'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.34%  '
$ws.Range('E27').Value = '  +4.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.28'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.24'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '28.63'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '640.98'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.41'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.25'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.107'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '56.77'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.77%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '36.58'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.04%  '
$ws.Range('E39').Value = '  +1.26%  '
$ws.Range('D40').Value = '0.0₃0734'
$ws.Range('E40').Value = '  +12.52%  '
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.62'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.124'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('D44').Value = '2.912.55'
$ws.Range('E44').Value = '  +1.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.94'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.69'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0396'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.61'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.05'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.22%  '
$ws.Range('E50').Value = '  +1.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '134.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.94%  '
